$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.3
$ws.Range("H2").Value = 2.6
$ws.Range("M2").Value = 1.84
$ws.Range("N2").Value = 2.02
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 1.67
$ws.Range("F3").Value = 2.1
$ws.Range("G3").Value = 2.8
$ws.Range("K3").Value = 1.47
$ws.Range("L3").Value = 2.5
$ws.Range("M3").Value = 2.52
$ws.Range("N3").Value = 1.46
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.6
$ws.Range("W3").Value = 1.5
$ws.Range("B4").Value = 45013.88541666666
$ws.Range("F4").Value = 1.6
$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 7
$ws.Range("K4").Value = 1.5
$ws.Range("L4").Value = 2.49
$ws.Range("M4").Value = 2.5
$ws.Range("N4").Value = 1.5
$ws.Range("V4").Value = 2.67
$ws.Range("W4").Value = 0.4
$ws.Range("F5").Value = 1.8
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 4.5
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 1.72
$ws.Range("V5").Value = 0.67
$ws.Range("W5").Value = 1.33
$ws.Range("AA5").Value = 1.24
$ws.Range("AB5").Value = 14.25
$ws.Range("AC5").Value = 4.6
$ws.Range("F6").Value = 2.11
$ws.Range("G6").Value = 3.15
$ws.Range("H6").Value = 3.5
$ws.Range("M6").Value = 2.25
$ws.Range("N6").Value = 1.57
$ws.Range("V6").Value = 1.75
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 3.3
$ws.Range("H7").Value = 2.15
$ws.Range("M7").Value = 1.75
$ws.Range("N7").Value = 1.97
$ws.Range("V7").Value = 1.62
$ws.Range("W7").Value = 1.57
$ws.Range("F8").Value = 3.15
$ws.Range("H8").Value = 2.13
$ws.Range("M8").Value = 2.39
$ws.Range("N8").Value = 1.51
$ws.Range("W8").Value = 0.57
$ws.Range("F9").Value = 2.5
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 2.6
$ws.Range("M9").Value = 1.7
$ws.Range("N9").Value = 2.1
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 2.38
$ws.Range("F10").Value = 1.82
$ws.Range("G10").Value = 3.4
$ws.Range("H10").Value = 3.6
$ws.Range("M10").Value = 1.92
$ws.Range("N10").Value = 1.88
$ws.Range("V10").Value = 2.15
$ws.Range("W10").Value = 1.5
$ws.Range("AA10").Value = 1.6
$ws.Range("AB10").Value = 9.6
$ws.Range("AC10").Value = 2.57
$ws.Range("AD10").Value = 1.23
$ws.Range("AE10").Value = 1.48
$ws.Range("AF10").Value = 1.73
$ws.Range("AG10").Value = 2.1
$ws.Range("AH10").Value = 2.7
$ws.Range("F11").Value = 2.74
$ws.Range("G11").Value = 3.35
$ws.Range("H11").Value = 2.06
$ws.Range("I11").Value = 1.03
$ws.Range("J11").Value = 14
$ws.Range("M11").Value = 1.75
$ws.Range("N11").Value = 1.95
$ws.Range("Q11").Value = 1.57
$ws.Range("R11").Value = 2.2
$ws.Range("V11").Value = 1.27
$ws.Range("W11").Value = 2.64
$ws.Range("G13").Value = 3.4
$ws.Range("H13").Value = 3.45
$ws.Range("M13").Value = 1.89
$ws.Range("N13").Value = 1.86
$ws.Range("AA13").Value = 1.58
$ws.Range("AB13").Value = 6
$ws.Range("AC13").Value = 2.65
$ws.Range("AD13").Value = 1.19
$ws.Range("AE13").Value = 1.36
$ws.Range("AF13").Value = 1.56
$ws.Range("AG13").Value = 1.9
$ws.Range("AH13").Value = 2.34
$ws.Range("F14").Value = 4.9
$ws.Range("G14").Value = 3.5
$ws.Range("H14").Value = 1.65
$ws.Range("M14").Value = 2.23
$ws.Range("N14").Value = 1.61
$ws.Range("F15").Value = 2.25
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 1.09
$ws.Range("J15").Value = 7.5
$ws.Range("K15").Value = 1.46
$ws.Range("L15").Value = 2.54
$ws.Range("M15").Value = 2.4
$ws.Range("N15").Value = 1.53
$ws.Range("O15").Value = 1.5
$ws.Range("P15").Value = 2.24
$ws.Range("Q15").Value = 2.15
$ws.Range("R15").Value = 1.62
$ws.Range("S15").Value = 1.32
$ws.Range("T15").Value = 1.33
$ws.Range("U15").Value = 1.5
$ws.Range("F16").Value = 2.9
$ws.Range("G16").Value = 2.7
$ws.Range("H16").Value = 2.5
$ws.Range("I16").Value = 1.11
$ws.Range("J16").Value = 6.5
$ws.Range("K16").Value = 1.57
$ws.Range("L16").Value = 2.2
$ws.Range("M16").Value = 2.85
$ws.Range("N16").Value = 1.37
$ws.Range("O16").Value = 1.6
$ws.Range("P16").Value = 2.2
$ws.Range("Q16").Value = 2.4
$ws.Range("R16").Value = 1.52
$ws.Range("S16").Value = 1.42
$ws.Range("T16").Value = 1.36
$ws.Range("U16").Value = 1.38
$ws.Range("F17").Value = 2.7
$ws.Range("G17").Value = 2.85
$ws.Range("H17").Value = 2.6
$ws.Range("K17").Value = 1.5
$ws.Range("L17").Value = 2.45
$ws.Range("M17").Value = 2.45
$ws.Range("N17").Value = 1.52
$ws.Range("AA17").Value = 1.86
$ws.Range("AB17").Value = 6.85
$ws.Range("AC17").Value = 2.3
$ws.Range("AD17").Value = 1.55
$ws.Range("AE17").Value = 1.93
$ws.Range("AF17").Value = 2.44
$ws.Range("F18").Value = 1.91
$ws.Range("G18").Value = 2.9
$ws.Range("H18").Value = 4.4
$ws.Range("I18").Value = 1.13
$ws.Range("J18").Value = 5.95
$ws.Range("K18").Value = 1.57
$ws.Range("L18").Value = 2.15
$ws.Range("O18").Value = 1.65
$ws.Range("P18").Value = 2.1
$ws.Range("Q18").Value = 2.63
$ws.Range("R18").Value = 1.44
$ws.Range("S18").Value = 1.17
$ws.Range("T18").Value = 1.32
$ws.Range("U18").Value = 1.95
$ws.Range("AF18").Value = 1.83
$ws.Range("AG18").Value = 1.86
$ws.Range("AH18").Value = 2.43
$ws.Range("F19").Value = 3.25
$ws.Range("G19").Value = 2.95
$ws.Range("H19").Value = 2.15
$ws.Range("I19").Value = 1.08
$ws.Range("J19").Value = 8.199999999999999
$ws.Range("K19").Value = 1.43
$ws.Range("L19").Value = 2.64
$ws.Range("M19").Value = 2.25
$ws.Range("N19").Value = 1.53
$ws.Range("O19").Value = 1.55
$ws.Range("P19").Value = 2.3
$ws.Range("Q19").Value = 2.05
$ws.Range("R19").Value = 1.7
$ws.Range("S19").Value = 1.65
$ws.Range("T19").Value = 1.32
$ws.Range("U19").Value = 1.3
$ws.Range("AF19").Value = 2.1
$ws.Range("F20").Value = 2.45
$ws.Range("G20").Value = 2.7
$ws.Range("H20").Value = 3.1
$ws.Range("I20").Value = 1.14
$ws.Range("J20").Value = 5.55
$ws.Range("K20").Value = 1.6
$ws.Range("L20").Value = 2.1
$ws.Range("M20").Value = 2.9
$ws.Range("N20").Value = 1.32
$ws.Range("O20").Value = 1.69
$ws.Range("P20").Value = 2.09
$ws.Range("Q20").Value = 2.2
$ws.Range("R20").Value = 1.62
$ws.Range("S20").Value = 1.35
$ws.Range("T20").Value = 1.4
$ws.Range("U20").Value = 1.42
$ws.Range("AE20").Value = 1.46
$ws.Range("AF20").Value = 2.25
$ws.Range("AG20").Value = 2.35
$ws.Range("AH20").Value = 3.25
